$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before the
#     existing "Late" column (old column N), pushing Late/heading/Outstanding
#     one column to the right (N->O, O->P, P->Q). ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()

# Match the inserted column's width to its neighbour (column M) as closely
# as the COM ColumnWidth property allows.
$ws.Columns("N").ColumnWidth = 9.8

# --- Make "Repayment schedule" the active/selected sheet & cell,
#     which also clears the previous tab selection on "Summary". ---
$ws.Activate()
$ws.Range("L14").Select()
